# Upload new version with timestamp
# Updates the "balance" ratio column (H) for several products in the
# day-sale shortage report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CONTAFEVER N 200MG/5ML SUSP. 120ML  (row 13): 19:0 -> 20:0
$ws.Range("H13").Value = "20:0"

# DIASMECT 20% SUSP. 60ML  (row 17): 15:0 -> 16:0
$ws.Range("H17").Value = "16:0"

# DOLIPRANE 1 GM 15 TABS.  (row 19): 6:1 -> 6:2
$ws.Range("H19").Value = "6:2"

# GASTROMOTIL 1MG/ML ORAL SUSP. 200ML  (row 20): 1:0 -> 2:0
$ws.Range("H20").Value = "2:0"

# MOTILIUM 10MG 40 F.C.TAB.  (row 22): 0:3 -> 1:0
$ws.Range("H22").Value = "1:0"

# STREPTOQUIN 20 TABLETS  (row 27): 4:0 -> 4:1
$ws.Range("H27").Value = "4:1"
